$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Class Statistics summary numbers (K3:L10 block)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 209
$ws.Range("L7").Value = 19

# Percent values stored as literal text (not numeric %), so force text format
# before assigning, otherwise Excel auto-converts "78.9%" into 0.789 with a
# percentage number format.
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "78.9%"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "75.4%"

# ---------------------------------------------------------------------------
# 2. Group Statistics block (rows 17 & 18)
# ---------------------------------------------------------------------------
$ws.Range("O17").Value = 37
$ws.Range("P17").Value = 5
$ws.Range("R17").NumberFormat = "@"
$ws.Range("R17").Value = "69.8%"
$ws.Range("S17").NumberFormat = "@"
$ws.Range("S17").Value = "79.9%"

$ws.Range("O18").Value = 46
$ws.Range("P18").Value = 1
$ws.Range("R18").NumberFormat = "@"
$ws.Range("R18").Value = "88.5%"
$ws.Range("S18").NumberFormat = "@"
$ws.Range("S18").Value = "76.0%"

# ---------------------------------------------------------------------------
# 3. "Recorded By" (column G) email-order swaps
# ---------------------------------------------------------------------------
$ws.Range("G24").Value = "haderreda2919@gmail.com, emp17.nada.h.attia@gmail.com"
$ws.Range("G180").Value = "dr.mohabelsawy@gmail.com, emp17.nada.h.attia@gmail.com"
$ws.Range("G235").Value = "dr.mohabelsawy@gmail.com, emp17.nada.h.attia@gmail.com"
$ws.Range("G213").Value = "fatma_shoukry@hotmail.com, drmohamedramadan50@gmail.com"

$swapRows = @(30,31,32,35,36,37,83,85,86,88,89,90,91,92,137,189,190,191,192,193,194,196,197,243,245,248,250,251)
foreach ($r in $swapRows) {
    $ws.Range("G$r").Value = "emp17.mariam.m.goda@gmail.com, 160534@med.asu.edu.eg"
}

# ---------------------------------------------------------------------------
# 4. Sessions that flip from "Not Recorded" (pink) to "Recorded" (green):
#    rows 150, 151, 155, 182 - copy the green "Recorded" formatting from a
#    known recorded row (row 2) onto A:I of each target row, then fill in
#    the recorded-by / students / status values.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A150:I150").PasteSpecial(-4122)
$ws.Range("A151:I151").PasteSpecial(-4122)
$ws.Range("A155:I155").PasteSpecial(-4122)
$ws.Range("A182:I182").PasteSpecial(-4122)

$ws.Range("G150").Value = "160715@med.asu.edu.eg"
$ws.Range("H150").Value = "29/64"
$ws.Range("I150").Value = "Recorded"

$ws.Range("G151").Value = "emp17.sara.h.ibrahim@gmail.com"
$ws.Range("H151").Value = "60/64"
$ws.Range("I151").Value = "Recorded"

$ws.Range("G155").Value = "160715@med.asu.edu.eg"
$ws.Range("H155").Value = "27/64"
$ws.Range("I155").Value = "Recorded"

$ws.Range("G182").Value = "Samarmaged9966@gmail.com"
$ws.Range("H182").Value = "34/62"
$ws.Range("I182").Value = "Recorded"
